$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (chars) to match the new, longer activity descriptions.
$ws.Columns("A").ColumnWidth = 84.42578125

# Row 3 ("Analizė" / 16) -> cleared, but keeps its bordered formatting.
$ws.Range("A3:B3").ClearContents()

# Rows 4-5 ("Projektavimas" / 40, "Realizavimas - kurti filma..." / 13)
# -> fully cleared (content + formatting), so they drop out of the sheet entirely.
$ws.Range("A4:B5").Clear()

# New block of rows 6-10: pasted-in content that kept the "no border" style.
$ws.Range("A6:B10").Borders.LineStyle = 0

$ws.Range("A6").Value = "Funkciniai reikalavimai kino studijos posistemei - klasių diagrama, panaudos atvejų diagrama"
$ws.Range("B6").Value = 6

$ws.Range("A7").Value = 'Technologijų analizė - "React" ir "React-Bootstrap"'
$ws.Range("B7").Value = 8

$ws.Range("A8").Value = 'Technologijų analizė - "Redux"'
$ws.Range("B8").Value = 6

$ws.Range("A9").Value = 'Technologijų analizė - "ASP.NET Core"'
$ws.Range("B9").Value = 8

$ws.Range("A10").Value = 'Technologijų analizė - "Entity Framework Core"'
$ws.Range("B10").Value = 7

# Row 11 keeps the bordered style.
$ws.Range("A11:B11").Borders.LineStyle = 1
$ws.Range("A11").Value = "Projektavimas - kino studijos posistemės panaudos atvejų sekų diagramos"
$ws.Range("B11").Value = 10

# Row 12 is unbordered again.
$ws.Range("A12:B12").Borders.LineStyle = 0
$ws.Range("A12").Value = "Realizavimas - filmo informacijos puslapio realizavimas"
$ws.Range("B12").Value = 10

# Rows 13-26: bordered, new/rearranged activity rows.
$ws.Range("A13:B26").Borders.LineStyle = 1

$ws.Range("A13").Value = 'Realizavimas – PA "Kurti filmą"'
$ws.Range("B13").Value = 13

$ws.Range("A14").Value = 'Realizavimas – PA "Redaguoti filmą"'
$ws.Range("B14").Value = 7

$ws.Range("A15").Value = 'Realizavimas –PA "Šalinti filmą"'
$ws.Range("B15").Value = 4

$ws.Range("A16").Value = "Realizavimas – kurti darbo skelbimą realizavimas"
$ws.Range("B16").Value = 10

$ws.Range("A17").Value = "Realizavimas – šalinti darbo skelbimą realizavimas"
$ws.Range("B17").Value = 4

$ws.Range("A18").Value = "Realizavimas – kino studijos registracijos realizavimas"
$ws.Range("B18").Value = 4

$ws.Range("A19").Value = "Realizavimas – kino filmų peržiūros puslapio realizavimas"
$ws.Range("B19").Value = 5

$ws.Range("A20").Value = "Realizavimas – darbo skelbimų peržiūros puslapio realizavimas"
$ws.Range("B20").Value = 4

$ws.Range("A21").Value = "Realizavimas – pranešimų peržiūros puslapio realizavimas"
$ws.Range("B21").Value = 4

$ws.Range("A22").Value = "Realizavimas – kino studijų statistikos puslapio realizavimas"
$ws.Range("B22").Value = 4

$ws.Range("A23").Value = "Realizavimas – kino studijos filmų ataskaitos puslapio realizavimas"
$ws.Range("B23").Value = 4

$ws.Range("A24").Value = "Kino studijos posistemės testavimas ir klaidų taisymas"
$ws.Range("B24").Value = 8

$ws.Range("A25").Value = "Kino studijos vartotojo vadovo sudarymas"
$ws.Range("B25").Value = 2

$ws.Range("A26").Value = "Ataskaita"
$ws.Range("B26").Value = 3

# Row 27 ("Testavimas" / 16) -> cleared, but keeps its bordered formatting.
$ws.Range("A27:B27").ClearContents()

# Rows 28-30 -> fully cleared (content + formatting), dropping out of the sheet.
$ws.Range("A28:B30").Clear()

# Selection ends up on the newly-typed "Ataskaita" row.
$ws.Range("A26:B26").Select()
